$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: emotion label was "Happy" but the files are the DIS (Disgust) set -> fix to "Disgust"
$ws.Range("A3").Value = "Disgust"

# Row 6: emotion label was "Surprise" but the files are the NEU (Neutral) set -> fix to "Neutral"
$ws.Range("A6").Value = "Neutral"

# Row 7: emotion label was "Disgust" but the files are the HAP (Happy) set -> fix to "Happy"
$ws.Range("A7").Value = "Happy"

# Update the last active selection cell, as recorded after the edit
$ws.Range("B10").Select()
